$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Range("B2").Value = "Vincix_ChallengeLv01"

$ws.Range("A3").Value = "OrchestratorQueueFolder"
$ws.Range("B3").Value = "Shared"
$ws.Range("C3").Value = "Folder name. The value must match a folder defined in Orchestrator and queue specified as OrchestratorQueueName should be created in this folder. For classic folders leave the value field empty."

$ws.Range("A8").Value = "VincixUrl"
$ws.Range("B8").Value = "https://www.rpahackathon.co.uk/login"

$ws.Range("A10").Value = "VincixCredential"
$ws.Range("B10").Value = "VincixCredential"

$ws.Range("A12").Value = "VincixLevel1Url"
$ws.Range("B12").Value = "https://www.rpahackathon.co.uk/l1"

$ws.Range("A14").Value = "ProcessName"
$ws.Range("B14").Value = "chrome"

$ws.Rows.Item(998).Delete()
$ws.Rows.Item(997).Delete()

$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 30

$ws.Range("B6").Select()
